# cv123023a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had three header/section-label rows that carried no numeric
# data of their own (row 2's "unnamed: 1_level_1" placeholder column
# header, plus the "situação do domicílio" and "grandes regiões e
# unidades da federação" section-title rows). The fix:
#   - gives the first data column its real header ("total") instead of
#     the pandas "unnamed" placeholder
#   - removes the two empty section-title rows so every label row lines
#     up with its own data again (rows 38/39 disappear because
#     everything above shifts up by two)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 currently holds "unnamed: 1_level_1" -> should read "total"
$ws.Range("B2").Value = "total"

# Remove the two label-only rows that have no data next to them.
# Delete the lower-numbered row last so the higher row number is still
# valid when it's deleted.
$ws.Rows("8:8").Delete()
$ws.Rows("5:5").Delete()
